# "Added 2020 results + clean up"
# - Remove the (now unused) "Assumptions" helper sheet.
# - Protect the four remaining worksheets.
# - Leave the EXPORT sheet as the active tab.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Delete the "Assumptions" sheet - it only held a tiny helper constant
# that's no longer referenced anywhere else in the workbook.
$wb.Worksheets.Item("Assumptions").Delete()

# Turn sheet protection on for the sheets that remain.
$wb.Worksheets.Item("Election Results by State").Protect()
$wb.Worksheets.Item("Uncontested Races").Protect()
$wb.Worksheets.Item("Uncontested by State PIVOT").Protect()
$wb.Worksheets.Item("EXPORT").Protect()

# EXPORT becomes the active tab.
$wb.Worksheets.Item("EXPORT").Activate()
